# fix(publipostage): Try to solve Excel emoji problem
#
# The "statut" column (column A) used emoji characters to flag the status
# of a clinical trial (blue/red/green/orange book emoji). Those emoji
# caused problems in the publipostage (mail-merge) pipeline, so they are
# replaced with plain-text / simple-symbol equivalents:
#
#   📘 (blue book)   -> ⚠️
#   📕 (red book)    -> -3
#   📗 (green book)  -> ✅
#   📙 (orange book) -> +3
#
# "-3" and "+3" look like numbers to Excel's smart-entry parser, so those
# two cells are forced to stay plain text (Text number format while the
# value is entered, then the format is cleared again) instead of silently
# turning into the numeric values -3 / 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "📘" = "⚠️"
    "📕" = "-3"
    "📗" = "✅"
    "📙" = "+3"
}

# Replacement values that Excel would otherwise auto-convert to numbers.
$numericLooking = @{
    "-3" = $true
    "+3" = $true
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    $key = [string]$val
    if (-not $map.ContainsKey($key)) { continue }

    $newVal = $map[$key]

    if ($numericLooking.ContainsKey($newVal)) {
        # Force text entry so "-3"/"+3" aren't reinterpreted as numbers,
        # then drop the temporary Text format again.
        $cell.NumberFormat = "@"
        $cell.Value2 = $newVal
        $cell.ClearFormats()
    } else {
        $cell.Value2 = $newVal
    }
}
